# Fix duplicate data and correct reconciliation results
# - Correct the Delta_new value for trade T001 (row 2)
# - Replace the duplicated T006-T015 block (rows 17-26) with the
#   missing unique trades T016-T025, updating PV_new / Delta_new
#   values accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: correct Delta_new for T001
$ws.Range("C2").Value = 0.44

# Row 17 (was duplicate T006) -> T016
$ws.Range("A17").Value = "T016"
$ws.Range("B17").Value = 93000
$ws.Range("C17").Value = 0.46

# Row 18 (was duplicate T007) -> T017
$ws.Range("A18").Value = "T017"
$ws.Range("B18").Value = -78500
$ws.Range("C18").Value = -0.76

# Row 19 (was duplicate T008) -> T018
$ws.Range("A19").Value = "T018"
$ws.Range("B19").Value = 104500
$ws.Range("C19").Value = 0.57

# Row 20 (was duplicate T009) -> T019 (no PV/Delta, stays blank)
$ws.Range("A20").Value = "T019"

# Row 21 (was duplicate T010) -> T020
$ws.Range("A21").Value = "T020"
$ws.Range("B21").Value = 83000
$ws.Range("C21").Value = 0.42

# Row 22 (was duplicate T011) -> T021
$ws.Range("A22").Value = "T021"
$ws.Range("B22").Value = -86000
$ws.Range("C22").Value = -0.87

# Row 23 (was duplicate T012) -> T022
$ws.Range("A23").Value = "T022"
$ws.Range("B23").Value = 99000
$ws.Range("C23").Value = 0.53

# Row 24 (was duplicate T013) -> T023 (no PV/Delta, stays blank)
$ws.Range("A24").Value = "T023"

# Row 25 (was duplicate T014) -> T024
$ws.Range("A25").Value = "T024"
$ws.Range("B25").Value = 73000
$ws.Range("C25").Value = 0.4

# Row 26 (was duplicate T015) -> T025
$ws.Range("A26").Value = "T025"
$ws.Range("B26").Value = -93000
$ws.Range("C26").Value = -0.84
